# Add multi regression model column: Total_beds_per1000 (column F)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, matching style of existing header row (bold/centered style used by A1:E1)
$ws.Range("F1").Value = "Total_beds_per1000"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108

$values = @(
    3.06,
    2.17,
    1.99,
    3.16,
    1.84,
    1.91,
    2.09,
    2.2,
    2.55,
    2.34,
    1.85,
    1.91,
    2.49,
    2.71,
    2.85,
    3.19,
    3.17,
    3.25,
    2.56,
    1.82,
    2.25,
    2.43,
    2.43,
    3.99,
    3.01,
    3.46,
    3.38,
    2.08,
    2.04,
    2.36,
    1.79,
    2.52,
    2.08,
    4.31,
    2.73,
    2.81,
    1.64,
    2.73,
    2,
    2.19,
    4.82,
    2.77,
    2.28,
    1.8,
    2.08,
    2.1,
    1.58,
    3.52,
    2.08,
    3.35
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
